$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 278, pushing existing rows 278:291 down to 279:292
$ws.Rows.Item(278).EntireRow.Insert()

# Populate the newly inserted row 278 with the new data record
$ws.Cells.Item(278, 1).Value = 10
$ws.Cells.Item(278, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(278, 3).Value = "La Araucanía"
$ws.Cells.Item(278, 4).Value = 44509
$ws.Cells.Item(278, 5).Value = 9
$ws.Cells.Item(278, 6).Value = "Fruta"
$ws.Cells.Item(278, 7).Value = 100103
$ws.Cells.Item(278, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(278, 9).Value = 100103006
$ws.Cells.Item(278, 10).Value = "Nectarín"
$ws.Cells.Item(278, 11).Value = "Early Glo"
$ws.Cells.Item(278, 12).Value = "Primera"
$ws.Cells.Item(278, 13).Value = 90
$ws.Cells.Item(278, 14).Value = 30000
$ws.Cells.Item(278, 15).Value = 32000
$ws.Cells.Item(278, 16).Value = 30778
$ws.Cells.Item(278, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(278, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(278, 19).Value = 1710
$ws.Cells.Item(278, 20).Value = 18

# Match the date number-format style used by column D elsewhere in the sheet
$ws.Cells.Item(278, 4).NumberFormat = $ws.Cells.Item(279, 4).NumberFormat
